# Edit script for "Persona - Al.docx"
# Applies the changes described by the target diff using Word COM automation.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. "Widowed, 4 years" -> "Widowed, 8 years"
# ---------------------------------------------------------------------------
Replace-Text "Widowed, 4 years" "Widowed, 8 years"

# ---------------------------------------------------------------------------
# 2. Restructure the lifestyle / work-history paragraphs.
#    Original paragraphs 7 & 8 (1-based, Word Paragraphs collection) are
#    rewritten and re-split into two new paragraphs.
# ---------------------------------------------------------------------------

$cr = [char]13

$p7 = $d.Paragraphs(7)
$p8 = $d.Paragraphs(8)
$rng78 = $d.Range($p7.Range.Start, $p8.Range.End)

$newParaA = "Al does not smoke, drinks moderately (2-3 beers a week, some weeks no drinks). A recreational cyclist, Al also likes to swim. Does not follow sports. Has been active with his church and has taught in Bible classes for ages 8 – 12, for a number of years. After retirement, he also began to work with the same age (8 – 12) group with a neighbourhood refugee support organization."
$newParaB = "Worked 40+ years for an architectural/building firm.Initially worked as a mechanical engineer, but then as mechanical designer. He is in demand there as an occasional trainer for new employees in the mechanical design department. Because of his past design work, Al enjoys drawing and watercolour painting."

$rng78.Text = $newParaA + $cr + $newParaB

# Paragraph 9: "Consequences of stroke..." (unchanged text, kept as-is).

# Paragraph 10: rehabilitation-unit paragraph - replace final sentence.
Replace-Text "Al has moved to assisted living accommodation but the intent is to help him overcome remaining problems so that he can return home, to work and to church." "Al has moved to assisted living accommodation but he struggles with the results of the stroke in terms of being motivated to improve his health."

# ---------------------------------------------------------------------------
# 3. "Walk through rooms..." bullet - "would have difficulty" -> "has difficulty"
# ---------------------------------------------------------------------------
Replace-Text "Walk through rooms, to bathroom, to kitchen. At present, he would have difficulty doing this." "Walk through rooms, to bathroom, to kitchen. At present, he has difficulty doing this."

# ---------------------------------------------------------------------------
# 4. "At present, in addition to walking problems..." bullet - big rewrite.
# ---------------------------------------------------------------------------
Replace-Text "At present, in addition to walking problems, Al would be unable to grasp and control the home implements such as cooking utensils. " "At present, in addition to walking problems, Al is often unable to grasp and control many home implements such as cooking utensils and has to rely on a ready-meal delivery system and also unable to draw or paint because of the lack of hand control. He finds the process to make his hands do as he wants very frustrating. "

# ---------------------------------------------------------------------------
# 5. Medication bullet - rewrite.
# ---------------------------------------------------------------------------
Replace-Text "Take medication on time. Initial assessments have shown Al is depressed with his recent stroke and feels out of control over what matters to him. This in turn leads him to not act according to instructions, such as medication. He presently has four types of medication: two are taken only once a day, a third is taken twice a day and should be taken, one, when he first gets up, and two, when he is ready to go to sleep; the fourth medication needs more precise administration and it is this one that Al failed to take correctly which led to some minor deterioration." "Take medication on time. Initial assessments have shown Al is depressed with his recent stroke and feels out of control over what matters to him. This in turn leads him to not act according to instructions, such as on medication. He presently has several types of medicine that need to be taken at specific times but despite labels on containers, Al failed to take one medicine correctly which led to some minor deterioration."

# ---------------------------------------------------------------------------
# 6. Robot bullet - drop "voice" and change ending.
# ---------------------------------------------------------------------------
Replace-Text "A robot which could take voice instructions" "A robot which could take instructions"
Replace-Text "in his home once released there. Items" "in his home again. Items"

# ---------------------------------------------------------------------------
# 7. "If these four areas..." paragraph - rewrite.
# ---------------------------------------------------------------------------
Replace-Text "to regain independence, he could be released to return to live home. His employers are keen to keep him in their employ and are making arrangements to give him ground-level office and access to ground-level rooms. Both a return to living at home and return to work would make a positive difference to Al, who misses deeply feeling like a functioning adult with much to give to others." "to regain independence, it could make a positive difference to Al, who deeply misses feeling like a functioning adult with much to give to others."

# ---------------------------------------------------------------------------
# 8. "Given that Al has worked..." paragraph - insert "past".
# ---------------------------------------------------------------------------
Replace-Text "It would tie in with his employment and give him an element of control" "It would tie in with his past employment and give him an element of control"
